# Generate Report for Handback
#
# This script reproduces a "handback" report generation pass: the
# localization pipeline finished writing the target + handback files for
# the two e2e documents in each language sheet, so the Overview/zh-cn/de-de
# "Status" columns flip from "Ready for handoff" to
# "Handed back: in sync with en-US", and the per-row Latest Target
# File / Latest Handback File / Latest Handback DateTime cells get filled
# in (with a genuine hyperlink on the new target-file cell, mirroring the
# existing handoff-file hyperlink).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$fileMdUrl = @{
    "744d9184-2d1d-4c6b-ba00-9afeef0cb72b.md" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38c10d60201e6daf7095f9c011b708d1d69e6541/e2e/744d9184-2d1d-4c6b-ba00-9afeef0cb72b.md"
    "d38922ec-6c63-40ba-996d-0fe14ec13819.md" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38c10d60201e6daf7095f9c011b708d1d69e6541/e2e/d38922ec-6c63-40ba-996d-0fe14ec13819.md"
}

# ---------------------------------------------------------------------
# Overview sheet: Status columns (zh-cn = E, de-de = F) for both rows.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Overview zh-cn / de-de columns got wider to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): Status column C, plus the
# Latest Target File (I) / Latest Handback File (J) /
# Latest Handback DateTime (K) columns for both data rows.
# ---------------------------------------------------------------------
$langSheets = @(
    @{
        Name = "zh-cn"
        HandbackDateTime = "2016-11-09 02:02:53"
        Row2 = @{
            File = "744d9184-2d1d-4c6b-ba00-9afeef0cb72b.md"
            Handback = "744d9184-2d1d-4c6b-ba00-9afeef0cb72b.ea76559ddb309da3e4e391a241a795cbf85f821a.zh-cn.xlf"
        }
        Row3 = @{
            File = "d38922ec-6c63-40ba-996d-0fe14ec13819.md"
            Handback = "d38922ec-6c63-40ba-996d-0fe14ec13819.113ba04d166113448b6836ea2b934ed91e02e3b1.zh-cn.xlf"
        }
    },
    @{
        Name = "de-de"
        HandbackDateTime = "2016-11-09 02:03:13"
        Row2 = @{
            File = "744d9184-2d1d-4c6b-ba00-9afeef0cb72b.md"
            Handback = "744d9184-2d1d-4c6b-ba00-9afeef0cb72b.ea76559ddb309da3e4e391a241a795cbf85f821a.de-de.xlf"
        }
        Row3 = @{
            File = "d38922ec-6c63-40ba-996d-0fe14ec13819.md"
            Handback = "d38922ec-6c63-40ba-996d-0fe14ec13819.113ba04d166113448b6836ea2b934ed91e02e3b1.de-de.xlf"
        }
    }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (C) for both rows.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Widen Status (C) and the now-populated Latest Target/Handback File
    # columns (I, J) so the longer values are readable.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    # Re-create the hyperlinks collection so the new "Latest Target File"
    # links (I2/I3) land interleaved with the existing handoff-file links
    # (A2/A3), in document order.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $fileMdUrl[$lang.Row2.File], "", "", $lang.Row2.File)
    $ws.Hyperlinks.Add($ws.Range("I2"), $fileMdUrl[$lang.Row2.File], "", "", $lang.Row2.File)

    $ws.Hyperlinks.Add($ws.Range("A3"), $fileMdUrl[$lang.Row3.File], "", "", $lang.Row3.File)
    $ws.Hyperlinks.Add($ws.Range("I3"), $fileMdUrl[$lang.Row3.File], "", "", $lang.Row3.File)

    # Latest Handback File (J) and Latest Handback DateTime (K).
    $ws.Range("J2").Value = $lang.Row2.Handback
    $ws.Range("K2").Value = $lang.HandbackDateTime

    $ws.Range("J3").Value = $lang.Row3.Handback
    $ws.Range("K3").Value = $lang.HandbackDateTime
}
